# Realestate Update resale numbers 2023-06-21 09:22
# Append a new data row (row 65) to the CityResaleNum sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 65

# Text columns (A-D) must stay as literal text, not get auto-converted
# to dates/numbers by Excel's type inference. Force the "Text" number
# format before assigning, then clear the formatting afterwards so the
# new row doesn't end up with a stray style applied.
$textRange = $ws.Range("A$row`:D$row")
$textRange.NumberFormat = "@"

$ws.Cells.Item($row, 1).Value = "2023-06-21"
$ws.Cells.Item($row, 2).Value = "09:22:15"
$ws.Cells.Item($row, 3).Value = "Wednesday"
$ws.Cells.Item($row, 4).Value = "25"

$textRange.ClearFormats()

# Numeric columns (E-T)
$ws.Cells.Item($row, 5).Value = 122158
$ws.Cells.Item($row, 6).Value = 133696
$ws.Cells.Item($row, 7).Value = 161958
$ws.Cells.Item($row, 8).Value = 133335
$ws.Cells.Item($row, 9).Value = 177304
$ws.Cells.Item($row, 10).Value = 114406
$ws.Cells.Item($row, 11).Value = 201536
$ws.Cells.Item($row, 12).Value = 225333
$ws.Cells.Item($row, 13).Value = 175508
$ws.Cells.Item($row, 14).Value = 103926
$ws.Cells.Item($row, 15).Value = 39283
$ws.Cells.Item($row, 16).Value = 33873
$ws.Cells.Item($row, 17).Value = 51940
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 35817
$ws.Cells.Item($row, 20).Value = -1
